# refactor: Enhance reusability of code
#
# Pull the "write a new UserName value into column C while preserving the
# cell's existing look (Text-formatted style id 2)" logic into one small,
# reusable helper instead of repeating the same steps for every row.
#
# Target state (per the source diff): the "UserName" column (C) for both
# data rows (2 and 3) gets a new value, 2521521663 (was 6651350157).
# The cells must stay genuine numeric cells (t="n") under their original
# style -- not be coerced into shared-string text, and not pick up a brand
# new style entry just because the format was nudged aside temporarily.

function Set-NumericValuePreservingStyle($Worksheet, $CellRef, $StyleDonorRef, $Value) {
    $target = $Worksheet.Range($CellRef)
    $styleDonor = $Worksheet.Range($StyleDonorRef)

    # The column is formatted as Text ("@"), so writing straight into it
    # would coerce the number into a shared string. Drop to the default
    # style first so the value commits as a real number...
    $target.Style = "Normal"
    $target.Value = $Value

    # ...then restore the original direct formatting by copying it over
    # from a cell that already carries the same style. This reuses the
    # existing cellXfs entry instead of minting a new, unused one.
    $styleDonor.Copy()
    $target.PasteSpecial(-4122)
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newUserName = 2521521663

Set-NumericValuePreservingStyle $ws "C2" "D2" $newUserName
Set-NumericValuePreservingStyle $ws "C3" "D3" $newUserName

$excel.CutCopyMode = $false
